$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the three obsolete "stability" metric columns (G:I) -
#    BankStability_score, ChannelStability_score, Stability_Mean.
#    Deleting the entire columns shifts everything from column J onward
#    three places to the left (J->G, K->H, ... W->T) and keeps the shared
#    string table / formulas in sync automatically.
# ---------------------------------------------------------------------------
$ws.Range("G1:I1").EntireColumn.Delete()

# ---------------------------------------------------------------------------
# 2. Refresh the remaining data with the newly supplied values (the commit
#    brought in new remote-sensing derived metrics, so several existing
#    columns now carry different numbers / need to be blanked out, and the
#    Riparian_Mean column now evaluates to a #NUM! error for every row).
# ---------------------------------------------------------------------------

# Row 2 - Ninemile 16-5
$ws.Range("A2").Value = "Ninemile 16-5"
$ws.Range("B2").Value = "Okanogan"
$ws.Range("C2").Value = "Ninemile Creek DS"
$ws.Range("E2").Value = "yes"
$ws.Range("G2").Value = 3
$ws.Range("H2").Value = 1
$ws.Range("M2").ClearContents()

# Row 3 - Salmon 16-11
$ws.Range("A3").Value = "Salmon 16-11"
$ws.Range("B3").Value = "Okanogan"
$ws.Range("C3").Value = "Salmon Creek-Lower"
$ws.Range("E3").Value = "yes"
$ws.Range("G3").Value = 3
$ws.Range("H3").Value = 5
$ws.Range("K3").ClearContents()
$ws.Range("M3").ClearContents()

# Row 4 - Salmon 16-6
$ws.Range("A4").Value = "Salmon 16-6"
$ws.Range("B4").Value = "Okanogan"
$ws.Range("C4").Value = "Salmon Creek-Lower"
$ws.Range("E4").Value = "yes"
$ws.Range("G4").Value = 3
$ws.Range("H4").Value = 5
$ws.Range("K4").ClearContents()
$ws.Range("L4").Value = 5
$ws.Range("M4").ClearContents()

# Row 5 - Salmon 16-9
$ws.Range("A5").Value = "Salmon 16-9"
$ws.Range("B5").Value = "Okanogan"
$ws.Range("C5").Value = "Salmon Creek-Lower"
$ws.Range("E5").Value = "yes"
$ws.Range("G5").Value = 3
$ws.Range("K5").ClearContents()
$ws.Range("M5").ClearContents()

# Row 6 - Tonasket 16-2
$ws.Range("A6").Value = "Tonasket 16-2"
$ws.Range("B6").Value = "Okanogan"
$ws.Range("C6").Value = "Tonasket Creek DS"
$ws.Range("E6").Value = "yes"
$ws.Range("G6").Value = 3
$ws.Range("H6").Value = 3
$ws.Range("K6").ClearContents()
$ws.Range("M6").ClearContents()

# ---------------------------------------------------------------------------
# 3. Riparian_Mean (column O) now errors out with #NUM! for every data row.
#    Use a worksheet-function call that genuinely raises #NUM!, then paste
#    the result back as a static value so no formula is stored in the file.
# ---------------------------------------------------------------------------
foreach ($r in 2..6) {
    $cell = $ws.Range("O$r")
    $cell.Value = $excel.WorksheetFunction.Sqrt(-1)
}

Write-Output "done"
